# 自动更新Excel文件 - 2026-01-11 23:13:18
#
# For every data row, recompute the "剩余" (days-remaining, column E) from
# the "总天" (total days, column D) and "开始时间" (start date, column F,
# stored as a plain YYYYMMDD integer) relative to "today". If the
# countdown has run out (remaining < 1), the cycle restarts: the start
# date is reset to "today" and the remaining count goes back to the full
# total days.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# "Today" for this run (the commit was produced the following day).
$today = Get-Date -Year 2026 -Month 1 -Day 12

function Get-DaysBetween($laterDate, $earlierDate) {
    # Subtracting DateTime objects directly isn't reliable here, but
    # OLE-Automation day numbers (whole days + time-of-day fraction) are -
    # since both values share the same time-of-day component it cancels
    # out cleanly, leaving the integer day count.
    $diffDays = $laterDate.ToOADate() - $earlierDate.ToOADate()
    return [Math]::Round($diffDays)
}

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count

for ($row = 2; $row -le $lastRow; $row++) {
    $totalDays = $ws.Cells.Item($row, 4).Value2
    $startRaw  = $ws.Cells.Item($row, 6).Value2

    if ($null -eq $totalDays -or $null -eq $startRaw) {
        continue
    }

    $startText = [string]$startRaw
    if ($startText.Length -ne 8) {
        # Malformed date (e.g. a typo like "202510929") - leave row as-is.
        continue
    }

    $year  = [int]$startText.Substring(0, 4)
    $month = [int]$startText.Substring(4, 2)
    $day   = [int]$startText.Substring(6, 2)

    if ($month -lt 1 -or $month -gt 12 -or $day -lt 1 -or $day -gt 31) {
        continue
    }

    $startDate = Get-Date -Year $year -Month $month -Day $day

    $elapsed = Get-DaysBetween $today $startDate
    $remaining = $totalDays - $elapsed

    if ($remaining -lt 1) {
        # Countdown expired - restart the cycle as of today.
        $ws.Cells.Item($row, 5).Value = $totalDays
        $ws.Cells.Item($row, 6).Value = [int]$today.ToString("yyyyMMdd")
    } else {
        $ws.Cells.Item($row, 5).Value = $remaining
    }
}
